# chore: adapt column header formatting to respective input file names (#7)
# Rename the "_old"/"_new" suffixed column headers to "_FV2210"/"_FV2304",
# wrap the data range in a table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) suffixes -------------------------------
# Columns A-J were "<name>_old" -> "<name>_FV2210"
$leftHeaders = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
for ($i = 0; $i -lt $leftHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $leftHeaders[$i]
}

# Column K ("diff") is unchanged.

# Columns L-U were "<name>_new" -> "<name>_FV2304"
$rightHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
for ($i = 0; $i -lt $rightHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $rightHeaders[$i]
}

# --- 2. Turn the used range into a table (Table1) -------------------------
$tableRange = $ws.Range("A1:U74")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header renaming, table creation and freeze panes applied"
